# Add a new "2022-Q4" sheet (right after "总计") and a matching summary row
# in "总计", per commit "feat: add 2022-Q4 data".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert the new quarter sheet right after the "总计" (total) sheet, so
#    the tab order becomes: 总计, 2022-Q4, 2022-Q3, 2022-Q2, 2021-Q1, 2020-Q4
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item(1)
$q4 = $wb.Worksheets.Add($null, $totalSheet)
$q4.Name = "2022-Q4"

# Columns B..G hold text (fund code / name / scale / position / ratio /
# market value) even when the text looks numeric, so force text format
# before writing -- otherwise "090019" becomes 90019, "1.18" becomes a
# float, etc. Columns A (index) and H (rank) stay plain numbers.
$q4.Range("B1:G5").NumberFormat = "@"

# Header row (same layout as the other quarterly sheets).
$q4.Cells.Item(1, 2).Value = "基金代码"
$q4.Cells.Item(1, 3).Value = "基金名称"
$q4.Cells.Item(1, 4).Value = "基金规模"
$q4.Cells.Item(1, 5).Value = "股票总仓位"
$q4.Cells.Item(1, 6).Value = "仓位占比"
$q4.Cells.Item(1, 7).Value = "持有市值(亿元)"
$q4.Cells.Item(1, 8).Value = "仓位排名"

# Holdings data rows.
$q4.Cells.Item(2, 1).Value = 0
$q4.Cells.Item(2, 2).Value = "090019"
$q4.Cells.Item(2, 3).Value = "大成景恒混合A"
$q4.Cells.Item(2, 4).Value = "1.18"
$q4.Cells.Item(2, 5).Value = "93.72"
$q4.Cells.Item(2, 6).Value = "1.85"
$q4.Cells.Item(2, 7).Value = "0.0218"
$q4.Cells.Item(2, 8).Value = 2

$q4.Cells.Item(3, 1).Value = 1
$q4.Cells.Item(3, 2).Value = "006038"
$q4.Cells.Item(3, 3).Value = "大成景恒混合C"
$q4.Cells.Item(3, 4).Value = "0.89"
$q4.Cells.Item(3, 5).Value = "93.72"
$q4.Cells.Item(3, 6).Value = "1.85"
$q4.Cells.Item(3, 7).Value = "0.0165"
$q4.Cells.Item(3, 8).Value = 2

$q4.Cells.Item(4, 1).Value = 2
$q4.Cells.Item(4, 2).Value = "001068"
$q4.Cells.Item(4, 3).Value = "国新国证新锐灵活配置混合"
$q4.Cells.Item(4, 4).Value = "0.21"
$q4.Cells.Item(4, 5).Value = "75.37"
$q4.Cells.Item(4, 6).Value = "6.50"
$q4.Cells.Item(4, 7).Value = "0.0136"
$q4.Cells.Item(4, 8).Value = 3

$q4.Cells.Item(5, 1).Value = 3
$q4.Cells.Item(5, 2).Value = "001797"
$q4.Cells.Item(5, 3).Value = "国新国证新利灵活配置混合"
$q4.Cells.Item(5, 4).Value = "0.02"
$q4.Cells.Item(5, 5).Value = "81.37"
$q4.Cells.Item(5, 6).Value = "7.75"
$q4.Cells.Item(5, 7).Value = "0.0016"
$q4.Cells.Item(5, 8).Value = 1

# Header row should keep the bold/centered/bordered look used by the other
# sheets' header rows (and by the "序号" index column).
$q4.Range("B1:H1").Font.Bold = $true
$q4.Range("B1:H1").HorizontalAlignment = -4108
$q4.Range("B1:H1").VerticalAlignment = -4160
$q4.Range("B1:H1").Borders.LineStyle = 1
$q4.Range("A2:A5").Font.Bold = $true
$q4.Range("A2:A5").HorizontalAlignment = -4108
$q4.Range("A2:A5").VerticalAlignment = -4160
$q4.Range("A2:A5").Borders.LineStyle = 1

# ---------------------------------------------------------------------------
# 2. Insert a new row 2 in "总计" for the 2022-Q4 summary, pushing the
#    existing quarters (2022-Q3, 2022-Q2, 2021-Q1, 2020-Q4) down by one row.
# ---------------------------------------------------------------------------
$totalSheet.Rows.Item(2).Insert()

$totalSheet.Cells.Item(2, 1).Value = 0
$totalSheet.Cells.Item(2, 2).Value = "2022-Q4"
$totalSheet.Cells.Item(2, 3).Value = 4
$totalSheet.Cells.Item(2, 4).Value = 0.05

# Keep the "序号" (index) column's bold/centered/bordered look.
$totalSheet.Range("A2:A6").Font.Bold = $true
$totalSheet.Range("A2:A6").HorizontalAlignment = -4108
$totalSheet.Range("A2:A6").VerticalAlignment = -4160
$totalSheet.Range("A2:A6").Borders.LineStyle = 1

# Renumber the "序号" (index) column so it stays 0..4 after the insert.
$totalSheet.Cells.Item(3, 1).Value = 1
$totalSheet.Cells.Item(4, 1).Value = 2
$totalSheet.Cells.Item(5, 1).Value = 3
$totalSheet.Cells.Item(6, 1).Value = 4
